$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 2243454.062
$ws.Range("B3").Value = 1497.816
$ws.Range("B4").Value = 1272.794
